$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = -21.705
    12 = -21.657
    15 = -22.044
    27 = -21.427
    28 = -21.777
    31 = -21.838
    32 = -21.664
    36 = -20.339
    38 = -19.741
    46 = -21.781
    54 = -21.862
    55 = -22.135
    56 = -22.086
    67 = -21.495
    69 = -21.721
    72 = -21.445
    73 = -19.981
    83 = -21.975
    86 = -22.093
    91 = -21.522
    93 = -21.259
    99 = -20.54
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
